$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of rows 8-13 (A8:C13) while keeping their formatting/style
$ws.Range("A8:C13").ClearContents()

# Update the active selection to B9
$ws.Range("B9").Select()
